$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2522.5715
$ws.Range("J18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("N18").Value = -4568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3144.2307
$ws.Range("I94").Value = 3144.2307
$ws.Range("K94").Value = 3144.2307
$ws.Range("M94").Value = -2693.2307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 38463532
$ws.Range("J103").Value = 38463532
$ws.Range("L103").Value = 115390596
$ws.Range("N103").Value = -115391768

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2350.8462
$ws.Range("I132").Value = 2416.5
$ws.Range("K132").Value = 7249.5
$ws.Range("M132").Value = -4719.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2267.7754
$ws.Range("I137").Value = 2096.7307
$ws.Range("J137").Value = 2461.1304
$ws.Range("K137").Value = 6290.1921
$ws.Range("L137").Value = 7383.3912
$ws.Range("M137").Value = -3740.1921
$ws.Range("N137").Value = -12483.3912

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5011.886
$ws.Range("I138").Value = 2508.1052
$ws.Range("J138").Value = 6914.76
$ws.Range("K138").Value = 7524.3156
$ws.Range("L138").Value = 20744.28
$ws.Range("M138").Value = -2384.3156
$ws.Range("N138").Value = -31024.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 907.8823
$ws.Range("I2").Value = 824.28
$ws.Range("K2").Value = 824.28
$ws.Range("M2").Value = -711.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10207.363
$ws.Range("I32").Value = 10141.274
$ws.Range("K32").Value = 10141.274
$ws.Range("M32").Value = -9854.273999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5838997
$ws.Range("I61").Value = 6255856
$ws.Range("K61").Value = 6255856
$ws.Range("M61").Value = -6255644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 907.8823
$ws.Range("I116").Value = 824.28
$ws.Range("K116").Value = 824.28
$ws.Range("M116").Value = 1469.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1880.1321
$ws.Range("I122").Value = 1914.0435
$ws.Range("K122").Value = 5742.1305
$ws.Range("M122").Value = -3292.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4767699
$ws.Range("I132").Value = 5055.273
$ws.Range("J132").Value = 10006607
$ws.Range("K132").Value = 15165.819
$ws.Range("L132").Value = 30019821
$ws.Range("M132").Value = -12635.819
$ws.Range("N132").Value = -30024881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5838997
$ws.Range("I136").Value = 6255856
$ws.Range("K136").Value = 18767568
$ws.Range("M136").Value = -18765018

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 907.8823
$ws.Range("I3").Value = 824.28
$ws.Range("K3").Value = 824.28
$ws.Range("M3").Value = -710.28

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5677.8887
$ws.Range("I20").Value = 6680.087
$ws.Range("J20").Value = 3904.7693
$ws.Range("K20").Value = 6680.087
$ws.Range("L20").Value = 3904.7693
$ws.Range("M20").Value = -6433.087
$ws.Range("N20").Value = -4398.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3091.2
$ws.Range("I86").Value = 2369.5
$ws.Range("K86").Value = 2369.5
$ws.Range("M86").Value = -1246.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3091.2
$ws.Range("I89").Value = 2369.5
$ws.Range("K89").Value = 11847.5
$ws.Range("M89").Value = -6231.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10000710
$ws.Range("I134").Value = 770.2857
$ws.Range("J134").Value = 33333904
$ws.Range("K134").Value = 2310.8571
$ws.Range("L134").Value = 100001712
$ws.Range("M134").Value = 224.1428999999998
$ws.Range("N134").Value = -100006782

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 9252.25
$ws.Range("I32").Value = 9252.25
$ws.Range("K32").Value = 9252.25
$ws.Range("M32").Value = -8936.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4665.8887
$ws.Range("I58").Value = 2332.5
$ws.Range("J58").Value = 9332.666999999999
$ws.Range("K58").Value = 2332.5
$ws.Range("L58").Value = 9332.666999999999
$ws.Range("M58").Value = -2129.5
$ws.Range("N58").Value = -9738.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2554
$ws.Range("I122").Value = 2197.6667
$ws.Range("J122").Value = 2859.4285
$ws.Range("K122").Value = 6593.000100000001
$ws.Range("L122").Value = 8578.2855
$ws.Range("M122").Value = -4143.000100000001
$ws.Range("N122").Value = -13478.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2641.1428
$ws.Range("I134").Value = 2641.1428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7923.428400000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5388.428400000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4665.8887
$ws.Range("I136").Value = 2332.5
$ws.Range("J136").Value = 9332.666999999999
$ws.Range("K136").Value = 6997.5
$ws.Range("L136").Value = 27998.001
$ws.Range("M136").Value = -4447.5
$ws.Range("N136").Value = -33098.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 358.44446
$ws.Range("J92").Value = 396.57144
$ws.Range("L92").Value = 1189.71432
$ws.Range("N92").Value = -3685.71432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 23610.953
$ws.Range("I133").Value = 23893
$ws.Range("K133").Value = 71679
$ws.Range("M133").Value = -66619

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3328.5
$ws.Range("I102").Value = 2938.3333
$ws.Range("K102").Value = 2938.3333
$ws.Range("M102").Value = -1316.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 927676.1
$ws.Range("I113").Value = 1244.6923
$ws.Range("J113").Value = 2648191.5
$ws.Range("K113").Value = 1244.6923
$ws.Range("L113").Value = 2648191.5
$ws.Range("M113").Value = 925.3077000000001
$ws.Range("N113").Value = -2652531.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5306496.5
$ws.Range("I132").Value = 3766.1482
$ws.Range("J132").Value = 21214688
$ws.Range("K132").Value = 11298.4446
$ws.Range("L132").Value = 63644064
$ws.Range("M132").Value = -8768.444600000001
$ws.Range("N132").Value = -63649124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 45999.777
$ws.Range("J134").Value = 45999.777
$ws.Range("L134").Value = 137999.331
$ws.Range("N134").Value = -143069.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7852.054
$ws.Range("I7").Value = 7328.6924
$ws.Range("K7").Value = 7328.6924
$ws.Range("M7").Value = -7216.6924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12000882
$ws.Range("I22").Value = 18857830
$ws.Range("K22").Value = 18857830
$ws.Range("M22").Value = -18857535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 12000882
$ws.Range("I27").Value = 18857830
$ws.Range("K27").Value = 18857830
$ws.Range("M27").Value = -18857723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3791.7874
$ws.Range("I40").Value = 3312.4412
$ws.Range("J40").Value = 5045.4614
$ws.Range("K40").Value = 3312.4412
$ws.Range("L40").Value = 5045.4614
$ws.Range("M40").Value = -3176.4412
$ws.Range("N40").Value = -5317.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1027
$ws.Range("I46").Value = 876.125
$ws.Range("K46").Value = 876.125
$ws.Range("M46").Value = -688.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1091.7894
$ws.Range("I55").Value = 446
$ws.Range("J55").Value = 1673
$ws.Range("K55").Value = 446
$ws.Range("L55").Value = 1673
$ws.Range("M55").Value = -273
$ws.Range("N55").Value = -2019

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 35753704
$ws.Range("I100").Value = 5196.8
$ws.Range("J100").Value = 125124980
$ws.Range("K100").Value = 5196.8
$ws.Range("L100").Value = 125124980
$ws.Range("M100").Value = -4655.8
$ws.Range("N100").Value = -125126062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3504.9375
$ws.Range("I122").Value = 3412.4905
$ws.Range("J122").Value = 3950.3635
$ws.Range("K122").Value = 10237.4715
$ws.Range("L122").Value = 11851.0905
$ws.Range("M122").Value = -7787.4715
$ws.Range("N122").Value = -16751.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7852.054
$ws.Range("I126").Value = 7328.6924
$ws.Range("K126").Value = 21986.0772
$ws.Range("M126").Value = -19516.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4051.1292
$ws.Range("I132").Value = 2747.111
$ws.Range("K132").Value = 8241.332999999999
$ws.Range("M132").Value = -5711.332999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5435.6665
$ws.Range("I136").Value = 3509.8
$ws.Range("J136").Value = 6398.6
$ws.Range("K136").Value = 10529.4
$ws.Range("L136").Value = 19195.8
$ws.Range("M136").Value = -7979.400000000001
$ws.Range("N136").Value = -24295.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7849.1055
$ws.Range("J62").Value = 13335.286
$ws.Range("L62").Value = 13335.286
$ws.Range("N62").Value = -14583.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7849.1055
$ws.Range("J65").Value = 13335.286
$ws.Range("L65").Value = 66676.42999999999
$ws.Range("N65").Value = -72916.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6145.1816
$ws.Range("I96").Value = 4459.6
$ws.Range("J96").Value = 7549.8335
$ws.Range("K96").Value = 4459.6
$ws.Range("L96").Value = 7549.8335
$ws.Range("M96").Value = -3086.6
$ws.Range("N96").Value = -10295.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2887.625
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 3120.2
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 6240.4
$ws.Range("M100").Value = -4459
$ws.Range("N100").Value = -7322.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3600.9714
$ws.Range("I107").Value = 2367.12
$ws.Range("J107").Value = 6685.6
$ws.Range("K107").Value = 7101.36
$ws.Range("L107").Value = 20056.8
$ws.Range("M107").Value = -5181.36
$ws.Range("N107").Value = -23896.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1404.76
$ws.Range("I113").Value = 1495.579
$ws.Range("K113").Value = 4486.737
$ws.Range("M113").Value = -2316.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2806.8333
$ws.Range("I126").Value = 3581.1765
$ws.Range("J126").Value = 926.2857
$ws.Range("K126").Value = 10743.5295
$ws.Range("L126").Value = 2778.8571
$ws.Range("M126").Value = -8273.529500000001
$ws.Range("N126").Value = -7718.8571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 76000
$ws.Range("J133").Value = 76000
$ws.Range("L133").Value = 76000
$ws.Range("N133").Value = -86120
